$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 6 and 7 (the dataset now only has 4 data rows instead of 6)
$ws.Rows("6:7").Delete()

# Update data rows 2-5 with the recalculated TPM values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppa"
$ws.Range("C2").Value = "Npr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3030883333333333
$ws.Range("H2").Value = 0.909265
$ws.Range("I2").Value = 0.5850568929085261
$ws.Range("J2").Value = 0.585056892908526
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1644733333333333
$ws.Range("N2").Value = 0.49342
$ws.Range("O2").Value = 0.2118645341380481
$ws.Range("P2").Value = 0.2118645341380481
$ws.Range("Q2").Value = 0.04984994847777778
$ws.Range("R2").Value = 0.4486495363
$ws.Range("S2").Value = 0.1239528060603188
$ws.Range("T2").Value = 0.1239528060603188
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppa"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3030883333333333
$ws.Range("H3").Value = 0.909265
$ws.Range("I3").Value = 0.5850568929085261
$ws.Range("J3").Value = 0.585056892908526
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6118403333333333
$ws.Range("N3").Value = 1.835521
$ws.Range("O3").Value = 0.7881354658619518
$ws.Range("P3").Value = 0.7881354658619518
$ws.Range("Q3").Value = 0.1854416668961111
$ws.Range("R3").Value = 1.668975002065
$ws.Range("S3").Value = 0.4611040868482073
$ws.Range("T3").Value = 0.4611040868482072
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Nppa"
$ws.Range("C4").Value = "Npr3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.214961
$ws.Range("H4").Value = 0.644883
$ws.Range("I4").Value = 0.4149431070914739
$ws.Range("J4").Value = 0.4149431070914739
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1644733333333333
$ws.Range("N4").Value = 0.49342
$ws.Range("O4").Value = 0.2118645341380481
$ws.Range("P4").Value = 0.2118645341380481
$ws.Range("Q4").Value = 0.03535535220666666
$ws.Range("R4").Value = 0.31819816986
$ws.Range("S4").Value = 0.08791172807772933
$ws.Range("T4").Value = 0.08791172807772933
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nppa"
$ws.Range("C5").Value = "Npr3"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.214961
$ws.Range("H5").Value = 0.644883
$ws.Range("I5").Value = 0.4149431070914739
$ws.Range("J5").Value = 0.4149431070914739
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.6118403333333333
$ws.Range("N5").Value = 1.835521
$ws.Range("O5").Value = 0.7881354658619518
$ws.Range("P5").Value = 0.7881354658619518
$ws.Range("Q5").Value = 0.1315218098936667
$ws.Range("R5").Value = 1.183696289043
$ws.Range("S5").Value = 0.3270313790137445
$ws.Range("T5").Value = 0.3270313790137445